$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the candidate record on row 2 with new iAuthor test case values
$ws.Range("A2").Value = "Owtxt234"
$ws.Range("B2").Value = 231028225
$ws.Range("C2").Value = "emhhupi52"
$ws.Range("D2").Value = "JT$!v89d"
$ws.Range("F2").Value = "HHuampQR"
$ws.Range("G2").Value = "Ewph"
